$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 31608322.37
$ws.Range("P2").Value = 322855144.56
$ws.Range("Q2").Value = 285172890.42
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 189922787.32
$ws.Range("T2").Value = 189922787.32
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = 30026260.11
$ws.Range("W2").Value = 22022054.01
$ws.Range("X2").Value = 5545226.85
$ws.Range("Y2").Value = 34235679.06
$ws.Range("Z2").Value = 33660628.27
$ws.Range("AA2").Value = 2052305.9
$ws.Range("AG2").Value = 3818369.87
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = ""
$ws.Range("AS2").Value = 25008822.37
$ws.Range("AT2").Value = ""
